$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column F header
$ws.Range("F1").Value = "importance"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Importance values for rows 2..128 (corresponding to data rows)
$values = @(100,95,90,85,80,65,20,20,65,50,25,20,15,20,20,15,10,5,50,45,35,35,50,15,15,20,40,75,20,35,40,30,25,30,25,35,15,15,15,15,15,20,20,20,35,30,25,25,25,25,25,25,25,25,25,25,25,25,25,25,25,25,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,20,15,15,15,1,15,15,10,10,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 6).Value = $v
    $row++
}
